# Update DTT Test Hour Log: add record for "User story 4" (and the
# preceding "User story 3" row that was still blank), per commit
# "update record for user story 4".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 - User story 3
$ws.Range("A7").Value = "User story 3 completed"
$ws.Range("B7").Value = 5
$ws.Range("C7").Value = "7/21/2024"
$ws.Range("D7").Value = "Everything explained on story 3 implemented such as map, house detail, and url launchers."
$ws.Range("D7").WrapText = $true
$ws.Rows.Item(7).RowHeight = 35

# Row 8 - User story 4
$ws.Range("A8").Value = "User story 4 completed "
$ws.Range("B8").Value = 3
$ws.Range("C8").Value = "7/21/2024"
$ws.Range("D8").Value = "Made design, and link "

# Recalculate totals and move the active selection like in the saved file
$excel.Calculate()
[void]$ws.Range("A11").Select()
